# Add a new worksheet "test_transfer" right after "test", copy the
# existing data over, and append a new "inAppDeliveryCode" column.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("test")

# New sheet, inserted right after "test".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "test_transfer"

# --- Row 1: headers (copied from "test", plus the new column) ---
$ws2.Range("A1").Value2 = $ws1.Range("A1").Value2
$ws2.Range("B1").Value2 = $ws1.Range("B1").Value2
$ws2.Range("C1").Value2 = $ws1.Range("C1").Value2
$ws2.Range("C1").NumberFormat = "yyyy\-mm\-dd"
$ws2.Range("D1").Value2 = $ws1.Range("D1").Value2
$ws2.Range("E1").Value2 = $ws1.Range("E1").Value2
$ws2.Range("F1").Value2 = $ws1.Range("F1").Value2
$ws2.Range("G1").Value2 = $ws1.Range("G1").Value2
$ws2.Range("H1").Value2 = $ws1.Range("H1").Value2
$ws2.Range("I1").Value2 = $ws1.Range("I1").Value2
$ws2.Range("J1").Value2 = "inAppDeliveryCode"

# --- Row 2: data (copied from "test", plus the new column) ---
$ws2.Range("A2").Value2 = $ws1.Range("A2").Value2
$ws2.Range("B2").Value2 = $ws1.Range("B2").Value2
$ws2.Range("C2").Value2 = $ws1.Range("C2").Value2
$ws2.Range("C2").NumberFormat = "yyyy\-mm\-dd"
$ws2.Range("D2").Value2 = $ws1.Range("D2").Value2
$ws2.Range("E2").Value2 = $ws1.Range("E2").Value2
$ws2.Range("F2").Value2 = $ws1.Range("F2").Value2
$ws2.Range("G2").Value2 = $ws1.Range("G2").Value2
$ws2.Range("H2").Value2 = $ws1.Range("H2").Value2
$ws2.Range("I2").Value2 = $ws1.Range("I2").Value2
$ws2.Range("J2").Value2 = "Y8P8ECFN8"

# Match the date-style formatting used for the dateOfBirth column (C)
# on the new inAppDeliveryCode column's header + value cells.
$ws2.Range("J1").NumberFormat = "yyyy\-mm\-dd"
$ws2.Range("J2").NumberFormat = "yyyy\-mm\-dd"

# Selections matching the saved workbook: full range selected on the
# original sheet, column J selected (and active) on the new sheet.
[void]$ws1.Range("A1:I2").Select()
[void]$ws2.Columns("J:J").Select()
[void]$ws2.Activate()
